# auto:changes in the call log form
#
# This script reproduces, via Excel COM-interop, the edits made to the
# "call" XLSForm workbook (survey / choices / settings sheets).
#
# Strategy: operate on the "survey" sheet (sheet1) from the bottom of the
# sheet upward so that row numbers used below always refer to the correct
# (not-yet-shifted) rows at the time each statement runs.

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")
$choices = $wb.Worksheets.Item("choices")

# ---------------------------------------------------------------------
# SURVEY sheet edits (bottom-up so earlier row numbers stay valid)
# ---------------------------------------------------------------------

# 1) note_verified label: add "STOP " warning (row 52)
$survey.Range("C52").Value = "Patient was not verified, STOP please end the call."

# 2) Remove the old address-verification block of 4 rows:
#    verify / match / if_no / if_yes  (rows 47-50)
$survey.Range("A47:A50").EntireRow.Delete()

# 3) Append the little "fa-user" icon markup to the dob note (old row 45)
$survey.Range("C45").Value = '${tsis_ctx}<I class="fa fa-user"></i>'

# 4) Insert two new "note" rows right after the dob note (old row 45):
#    - currnt_address : ${cur_address}
#    - note_1 : multi line identity verification script
$survey.Range("A46:A46").EntireRow.Insert()
$survey.Range("A46").Value = "note"
$survey.Range("B46").Value = "currnt_address"
$survey.Range("C46").Value = '${cur_address}'

$survey.Range("A47:A47").EntireRow.Insert()
$survey.Range("A47").Value = "note"
$survey.Range("B47").Value = "note_1"
$survey.Range("C47").Value = "Verify the Client’s identity:`n    1. Can you please remind me of your last name?`n    2. Please tell me your DOB? `n    3. What is your current address?`n    4. What clinic did you last visit?`n"

# 5) Append the "fa-user" icon markup to the last_name note (old row 44)
$survey.Range("C44").Value = '${aka_ctx}  |  ${yr_date_of_birth_ctx} yr   ${at_ctx}<I class="fa fa-user"></i>'

# 6) Insert the new "cur_address" calculate row right after "patient_aka"
#    (old row 31), pulling from the contact's curr_address field
$survey.Range("A32:A32").EntireRow.Insert()
$survey.Range("A32").Value = "calculate"
$survey.Range("B32").Value = "cur_address"
$survey.Range("C32").Value = "NO_LABEL"
$survey.Range("G32").Value = "../inputs/contact/curr_address"

# 7) Insert the new hidden "curr_address" string field right after "sex"
#    (old row 15), inside the inputs/contact group
$survey.Range("A16:A16").EntireRow.Insert()
$survey.Range("A16").Value = "string"
$survey.Range("B16").Value = "curr_address"
$survey.Range("C16").Value = "NO_LABEL"
$survey.Range("E16").Value = "hidden"

# ---------------------------------------------------------------------
# CHOICES sheet edits
# ---------------------------------------------------------------------

# Remove the "added_lab_appointment" choice option from the actions list
$choices.Range("A5:A5").EntireRow.Delete()

Write-Host "done"
